$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.051999999999989
$ws.Range("B21").Value = 5.689699999999993
$ws.Range("B23").Value = 5.6693
$ws.Range("B25").Value = 5.920299999999994
